# Swap the data between row 9 and row 10 (columns A, B, E, F, G, H, Q, R, S)
# as described by the diff: the two rows exchange their identifying data while
# the remaining columns (D, P, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY)
# stay identical between the two rows, so no change is needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "S")

foreach ($col in $cols) {
    $addr9  = "$col" + "9"
    $addr10 = "$col" + "10"

    $val9  = $ws.Range($addr9).Value2
    $val10 = $ws.Range($addr10).Value2

    $ws.Range($addr9).Value2  = $val10
    $ws.Range($addr10).Value2 = $val9
}
